$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.576.23"

$ws.Range("D3").Value = "2.699.88"
$ws.Range("E3").Value = "  +2.11%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'598.97"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "'160.05"
$ws.Range("E6").Value = "  +2.13%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "2.699.28"
$ws.Range("E9").Value = "  +2.12%  "

$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.30"
$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.360"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").Value = "'28.31"
$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "68.539.92"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "2.694.85"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").Value = "'11.99"
$ws.Range("E19").Value = "  +5.60%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.71"
$ws.Range("E20").Value = "  +4.13%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'366.80"
$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("D22").Value = "'4.56"
$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("D23").Value = "'4.90"
$ws.Range("E23").Value = "  +1.77%  "

$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +3.21%  "

$ws.Range("D25").Value = "'74.49"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("E27").Value = "  +4.17%  "

$ws.Range("D28").Value = "2.836.36"

$ws.Range("D29").Value = "'0.0000105"
$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").Value = "'573.41"
$ws.Range("E31").Value = "  +3.36%  "

$ws.Range("D32").Value = "'8.27"
$ws.Range("E32").Value = "  +3.22%  "

$ws.Range("D33").Value = "'1.45"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").Value = "'1.95"
$ws.Range("E34").Value = "  +5.43%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.132"
$ws.Range("E35").Value = "  +3.07%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.65"
$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'20.02"
$ws.Range("E38").Value = "  +3.27%  "

$ws.Range("D39").Value = "'160.99"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "'0.381"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("E41").Value = "  +2.04%  "

$ws.Range("D42").Value = "'5.42"
$ws.Range("E42").Value = "  +1.77%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'17.86"
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.65"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0317"
$ws.Range("E45").Value = "  -6.39%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "'158.35"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").Value = "'3.95"
$ws.Range("E48").Value = "  +5.78%  "

$ws.Range("E49").Value = "  +5.01%  "

$ws.Range("D50").Value = "'0.601"
$ws.Range("E50").Value = "  +7.14%  "

$ws.Range("D51").Value = "'22.10"
$ws.Range("E51").Value = "  +0.19%  "
